$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 and 7 (data reduced from 6 to 4 data rows)
$ws.Rows("6:7").Delete()

# Update remaining rows 2-5 with new TPM-derived values
# Row 2
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1149353333333333
$ws.Range("H2").Value = 0.344806
$ws.Range("I2").Value = 0.05719122335670149
$ws.Range("J2").Value = 0.05719122335670149
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03801766666666666
$ws.Range("N2").Value = 0.114053
$ws.Range("O2").Value = 0.001128449675396954
$ws.Range("P2").Value = 0.001128449675396954
$ws.Range("Q2").Value = 0.004369573190888889
$ws.Range("R2").Value = 0.039326158718
$ws.Range("S2").Value = (6.453741743242446 / 100000)
$ws.Range("T2").Value = (6.453741743242446 / 100000)

# Row 3
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1149353333333333
$ws.Range("H3").Value = 0.344806
$ws.Range("I3").Value = 0.05719122335670149
$ws.Range("J3").Value = 0.05719122335670149
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 33.652157
$ws.Range("N3").Value = 100.956471
$ws.Range("O3").Value = 0.998871550324603
$ws.Range("P3").Value = 0.9988715503246031
$ws.Range("Q3").Value = 3.867821882180666
$ws.Range("R3").Value = 34.810396939626
$ws.Range("S3").Value = 0.05712668593926906
$ws.Range("T3").Value = 0.05712668593926906

# Row 4
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.894732
$ws.Range("H4").Value = 5.684196
$ws.Range("I4").Value = 0.9428087766432985
$ws.Range("J4").Value = 0.9428087766432984
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03801766666666666
$ws.Range("N4").Value = 0.114053
$ws.Range("O4").Value = 0.001128449675396954
$ws.Range("P4").Value = 0.001128449675396954
$ws.Range("Q4").Value = 0.07203328959866667
$ws.Range("R4").Value = 0.648299606388
$ws.Range("S4").Value = 0.001063912257964529
$ws.Range("T4").Value = 0.001063912257964529

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.894732
$ws.Range("H5").Value = 5.684196
$ws.Range("I5").Value = 0.9428087766432985
$ws.Range("J5").Value = 0.9428087766432984
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 33.652157
$ws.Range("N5").Value = 100.956471
$ws.Range("O5").Value = 0.998871550324603
$ws.Range("P5").Value = 0.9988715503246031
$ws.Range("Q5").Value = 63.76181873692399
$ws.Range("R5").Value = 573.8563686323159
$ws.Range("S5").Value = 0.9417448643853339
$ws.Range("T5").Value = 0.9417448643853339

